$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (B13/C13 = "519033 - Carlos Yujiro Shigue", no A13) is removed;
# this shifts all subsequent rows (and their heights) up by one.
$ws.Rows.Item(13).Delete()

# Row 10: Objetivos value changes to the docente text.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13 (was row 14): Programa resumido / Semestral
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 14 (was row 15): Short syllabus
$ws.Range("A14").Value = "Short syllabus:"

# Row 15 (was row 16): Programa / 01/01/2012
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

# Row 16 (was row 17): Syllabus
$ws.Range("A16").Value = "Syllabus:"

# Row 17 (was row 18): Avaliação
$ws.Range("A17").Value = "Avaliação:"

# Row 18 (was row 19): Método / docente text
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

# Row 19 (was row 20): Critério / Aulas expositivas...
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e práticas ministradas em laboratório."
$ws.Range("C19").Value = "Aulas expositivas e práticas ministradas em laboratório."

# Row 20 (was row 21): Norma de recuperação / Média ponderada...
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# Row 21 (was row 22, but already removed by the shift): Bibliografia / Média ponderada...
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C21").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
